$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.9999966894208786
$ws.Range("E2").Value = 0.9999966894208786

$ws.Range("D3").Value = 0.8551576537463798
$ws.Range("E3").Value = 0.8551576537463798

$ws.Range("D4").Value = [double]"7.953810386309141E-52"
$ws.Range("E4").Value = [double]"7.953810386309141E-52"

$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

$ws.Range("D6").Value = 0.9999954604946731
$ws.Range("E6").Value = 0.9999954604946731

$ws.Range("D8").Value = 0.9999999999999936
$ws.Range("E8").Value = [double]"6.439293542825908E-15"

$ws.Range("D10").Value = [double]"2.539196964916301E-27"

$ws.Range("D11").Value = [double]"6.091091740398034E-07"
$ws.Range("E11").Value = 0.999999390890826
$ws.Range("F11").Value = 118.2373275756836
